$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.100.04"
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "'1.630.38"
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'216.21"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'0.514"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  +0.74%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.253"
$ws.Range("D8").Style = $ws.Range("C8").Style
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").Value = "'20.14"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").Value = "'1.857.72"
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "'1.626.78"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  -1.15%  "
$ws.Range("D14").Value = "'4.11"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "'65.73"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "  -3.45%  "
$ws.Range("D17").Value = "'27.071.78"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.38"
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("B23").Value = "Toncoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D23").Value = "'2.50"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'9.11"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'147.20"
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = "  -0.73%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'7.39"
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("E28").Value = "  -1.34%  "
$ws.Range("D29").Value = "'15.59"
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").Value = "'0.0504"
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "'3.36"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").Value = "'3.00"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "'1.300.98"
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = "  +2.27%  "
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "'0.541"
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("D39").Value = "'0.843"
$ws.Range("D39").Style = $ws.Range("C39").Style
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'2.26"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  +2.12%  "
$ws.Range("D42").Value = "'0.807"
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("E43").Value = "  -2.13%  "
$ws.Range("D44").Value = "'1.766.95"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  -1.36%  "
$ws.Range("D45").Value = "'62.23"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = "  -1.47%  "
$ws.Range("D46").Value = "'90.58"
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = "  -2.22%  "
$ws.Range("D47").Value = "'1.60"
$ws.Range("D47").Style = $ws.Range("C47").Style
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").Value = "'0.801"
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = "  +19.26%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("D51").Value = "'7.56"
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = "  -2.05%  "
